$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - JESUS REMBERTO DEVIA MORALES
$ws.Range("C16").Value = "1143357464"
$ws.Range("D16").Value = "JESUS REMBERTO DEVIA MORALES"
$ws.Range("E16").Value = "1702"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 689456

# Row 17 - HERNANDO JOSE MARMOLEJO FLOREZ (period 1702)
$ws.Range("C17").Value = "73140025"
$ws.Range("D17").Value = "HERNANDO JOSE MARMOLEJO FLOREZ"
$ws.Range("E17").Value = "1702"
$ws.Range("F17").Value = 27578
$ws.Range("G17").Value = 689456

# Row 18 - HERNANDO JOSE MARMOLEJO FLOREZ (period 1703)
$ws.Range("C18").Value = "73140025"
$ws.Range("D18").Value = "HERNANDO JOSE MARMOLEJO FLOREZ"
$ws.Range("E18").Value = "1703"
$ws.Range("F18").Value = 27578
$ws.Range("G18").Value = 689456

# Row 19 - ORLANDO ENRIQUE TILVEZ MARRUGO (period 1801)
$ws.Range("C19").Value = "1128046739"
$ws.Range("D19").Value = "ORLANDO ENRIQUE TILVEZ MARRUGO"
$ws.Range("E19").Value = "1801"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 737717

# Row 20 - JENNYFER ALTAMAR MORALES (period 1710)
$ws.Range("C20").Value = "45544387"
$ws.Range("D20").Value = "JENNYFER ALTAMAR MORALES"
$ws.Range("E20").Value = "1710"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 737717

# Row 21 - ORLANDO ENRIQUE TILVEZ MARRUGO (period 1704)
$ws.Range("C21").Value = "1128046739"
$ws.Range("D21").Value = "ORLANDO ENRIQUE TILVEZ MARRUGO"
$ws.Range("E21").Value = "1704"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 737717

# Row 22 - ORLANDO ENRIQUE TILVEZ MARRUGO (period 1802)
$ws.Range("C22").Value = "1128046739"
$ws.Range("D22").Value = "ORLANDO ENRIQUE TILVEZ MARRUGO"
$ws.Range("E22").Value = "1802"
$ws.Range("F22").Value = 15738
$ws.Range("G22").Value = 737717
